$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, pushing existing rows 18-22 down to 19-23
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with data (same pattern as surrounding rows, new values for
# the varying columns: D, J, K, L, M, P)
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44988
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100114007
$ws.Range("G18").Value = "Jengibre"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17400
$ws.Range("N18").Value = "$/caja 13 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 1338
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "Hortaliza"
